$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (Tue Jan  2 20:15:30 UTC 2024 GitHub Actions run):
# updated Price (D) / Volume(1h) (E) figures for every coin, plus Maker
# jumping up the ranking (new row 44), which shifts Algorand, Celestia and
# FirstDigitalUSD down one row each (rows 45-47).
#
# Every value is written with a leading apostrophe so Excel always stores it
# as literal text (mirroring the original inlineStr cells) instead of auto-
# coercing number-looking strings (e.g. "8.50", "1.00", "0.110") to numbers,
# which would silently drop the significant trailing zeros.


# Row 2
$ws.Range("D2").Value = "'45.247.75"
$ws.Range("E2").Value = "'  +2.99%  "

# Row 3
$ws.Range("D3").Value = "'2.377.60"
$ws.Range("E3").Value = "'  +1.29%  "

# Row 4
$ws.Range("E4").Value = "'  -0.31%  "

# Row 5
$ws.Range("D5").Value = "'311.73"
$ws.Range("E5").Value = "'  -0.51%  "

# Row 6
$ws.Range("D6").Value = "'108.79"
$ws.Range("E6").Value = "'  -0.91%  "

# Row 7
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "'  -0.17%  "

# Row 8
$ws.Range("E8").Value = "'  -0.09%  "

# Row 9
$ws.Range("D9").Value = "'0.619"
$ws.Range("E9").Value = "'  -0.63%  "

# Row 10
$ws.Range("D10").Value = "'41.11"
$ws.Range("E10").Value = "'  -0.98%  "

# Row 11
$ws.Range("D11").Value = "'0.0919"
$ws.Range("E11").Value = "'  -0.45%  "

# Row 12
$ws.Range("D12").Value = "'8.50"
$ws.Range("E12").Value = "'  -1.78%  "

# Row 13
$ws.Range("D13").Value = "'0.110"
$ws.Range("E13").Value = "'  +1.22%  "

# Row 14
$ws.Range("D14").Value = "'0.981"
$ws.Range("E14").Value = "'  -3.64%  "

# Row 15
$ws.Range("D15").Value = "'2.738.48"
$ws.Range("E15").Value = "'  +1.50%  "

# Row 16
$ws.Range("D16").Value = "'15.32"
$ws.Range("E16").Value = "'  -1.77%  "

# Row 17
$ws.Range("D17").Value = "'2.368.54"
$ws.Range("E17").Value = "'  +1.12%  "

# Row 18
$ws.Range("D18").Value = "'45.223.83"
$ws.Range("E18").Value = "'  +3.12%  "

# Row 19
$ws.Range("D19").Value = "'14.72"
$ws.Range("E19").Value = "'  +12.17%  "

# Row 20
$ws.Range("D20").Value = "'7.30"
$ws.Range("E20").Value = "'  -4.49%  "

# Row 21
$ws.Range("E21").Value = "'  -0.90%  "

# Row 22
$ws.Range("D22").Value = "'73.38"
$ws.Range("E22").Value = "'  -1.39%  "

# Row 23
$ws.Range("D23").Value = "'3.50"
$ws.Range("E23").Value = "'  -0.64%  "

# Row 24
$ws.Range("D24").Value = "'260.21"
$ws.Range("E24").Value = "'  -3.61%  "

# Row 25
$ws.Range("D25").Value = "'2.30"
$ws.Range("E25").Value = "'  +0.30%  "

# Row 26
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "'  +0.62%  "

# Row 27
$ws.Range("D27").Value = "'11.17"
$ws.Range("E27").Value = "'  -0.33%  "

# Row 28
$ws.Range("D28").Value = "'7.27"
$ws.Range("E28").Value = "'  -4.99%  "

# Row 29
$ws.Range("E29").Value = "'  +1.96%  "

# Row 30
$ws.Range("D30").Value = "'0.0968"
$ws.Range("E30").Value = "'  +9.09%  "

# Row 31
$ws.Range("D31").Value = "'22.43"
$ws.Range("E31").Value = "'  -1.31%  "

# Row 32
$ws.Range("D32").Value = "'37.61"
$ws.Range("E32").Value = "'  -4.15%  "

# Row 33
$ws.Range("D33").Value = "'169.18"
$ws.Range("E33").Value = "'  +0.56%  "

# Row 34
$ws.Range("D34").Value = "'2.94"
$ws.Range("E34").Value = "'  +6.10%  "

# Row 35
$ws.Range("E35").Value = "'  -1.37%  "

# Row 36
$ws.Range("D36").Value = "'0.117"
$ws.Range("E36").Value = "'  +2.97%  "

# Row 37
$ws.Range("D37").Value = "'4.75"
$ws.Range("E37").Value = "'  -1.30%  "

# Row 38
$ws.Range("D38").Value = "'3.94"
$ws.Range("E38").Value = "'  +2.48%  "

# Row 39
$ws.Range("D39").Value = "'2.95"
$ws.Range("E39").Value = "'  +2.15%  "

# Row 40
$ws.Range("D40").Value = "'0.0355"
$ws.Range("E40").Value = "'  -3.55%  "

# Row 41
$ws.Range("D41").Value = "'1.76"
$ws.Range("E41").Value = "'  +1.98%  "

# Row 42
$ws.Range("D42").Value = "'99.94"
$ws.Range("E42").Value = "'  -4.62%  "

# Row 43
$ws.Range("D43").Value = "'69.71"
$ws.Range("E43").Value = "'  -3.17%  "

# Row 44
$ws.Range("B44").Value = "'Maker"
$ws.Range("C44").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'1.878.63"
$ws.Range("E44").Value = "'  +12.58%  "

# Row 45
$ws.Range("B45").Value = "'Algorand"
$ws.Range("C45").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.229"
$ws.Range("E45").Value = "'  -3.77%  "

# Row 46
$ws.Range("B46").Value = "'Celestia"
$ws.Range("C46").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "'12.97"
$ws.Range("E46").Value = "'  -3.06%  "

# Row 47
$ws.Range("B47").Value = "'FirstDigitalUSD"
$ws.Range("C47").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "'  -0.24%  "

# Row 48
$ws.Range("D48").Value = "'81.61"
$ws.Range("E48").Value = "'  +5.38%  "

# Row 49
$ws.Range("D49").Value = "'5.66"
$ws.Range("E49").Value = "'  +6.40%  "

# Row 50
$ws.Range("D50").Value = "'112.13"
$ws.Range("E50").Value = "'  -2.24%  "

# Row 51
$ws.Range("D51").Value = "'9.21"
$ws.Range("E51").Value = "'  +2.14%  "
